# Add a new test case row (fm21 / wip) to the "ftests" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- New data row 26 --------------------------------------------------
$ws.Range("B26").Value = "fm21"
$ws.Range("C26").Value = "Location deductibles with overall maximum policy deductible, and policy limit with IL back-allocation"
$ws.Range("D26").Value = "0,2"
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = "3,4,12"
$ws.Range("I26").Value = "wip"
$ws.Range("J26").Value = "wip"

# --- Match formatting of the surrounding rows --------------------------
$ws.Range("B25:C25").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D8:H8").Copy()
$ws.Range("D26:H26").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I25:J25").Copy()
$ws.Range("I26:J26").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Update view state to reflect the new selection --------------------
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("I27").Select()
